$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: Copy the "Meta description" paragraph (2nd paragraph, right
# after the title heading) - excluding its trailing paragraph mark so
# that pasting it elsewhere doesn't create a brand-new paragraph break
# (it will just splice its runs in place).
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaRange = $metaPara.Range
$metaCopyRange = $d.Range($metaRange.Start, $metaRange.End - 1)
$metaCopyRange.Select()
$word.Selection.Copy()

# ------------------------------------------------------------------
# Step 2: Remove that paragraph entirely from the top of the document.
# ------------------------------------------------------------------
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# Step 3: Paste the copied runs right before the final (italic) image
# prompt paragraph, at the very start of its text.
# ------------------------------------------------------------------
$total = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($total)
$lastRange = $lastPara.Range
$pasteStart = $lastRange.Start
$pastePoint = $d.Range($pasteStart, $pasteStart)
$pastePoint.Select()
$word.Selection.Paste()

# ------------------------------------------------------------------
# Step 4: Split the merged paragraph into two paragraphs again, right
# before the "Create a feature image" text, so the pasted meta runs
# form their own paragraph and the image-prompt text keeps its own.
# ------------------------------------------------------------------
$splitSearch = $d.Range($pasteStart, $d.Content.End)
$splitSearch.Find.Execute("Create a feature image") | Out-Null
$splitPoint = $d.Range($splitSearch.Start, $splitSearch.Start)
$splitPoint.InsertBefore("`r")

# ------------------------------------------------------------------
# Step 5: Rename the bold "Meta description" run (now near the end of
# the document) to the new heading text, and drop the leftover
# ": Read our review..." run that used to follow it (its content moves
# to the image-prompt paragraph below).
# ------------------------------------------------------------------
$renameRange = $d.Range($pasteStart, $d.Content.End)
$renameRange.Find.Execute("Meta description", $false, $false, $false, $false, $false, $true, 1, $false, `
    "Play Drift King Free: A Unique Game with High-speed Thrill", 2) | Out-Null

$dropRange = $d.Range($pasteStart, $d.Content.End)
$dropRange.Find.Execute(": Read our review of Drift King, a unique racing game with multiple bonuses, and play for free. Enjoy challenging gameplay and immerse in high-speed visuals.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ------------------------------------------------------------------
# Step 6: Replace the old image-generation prompt text with the
# meta-description copy text (keeping the paragraph's italic run).
# ------------------------------------------------------------------
$imgRange = $d.Range($pasteStart, $d.Content.End)
$imgRange.Find.Execute("Create a feature image for " + [char]34 + "Drift King" + [char]34 + " that showcases a happy warrior with a pair of glasses in cartoon style. The image should feature the game's signature yellow car and high-speed racing elements in the background. The warrior should be seen holding a trophy, with a confident and victorious expression on their face. The overall style should capture the fast-paced nature of the game and appeal to fans of both slot games and racing.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Read our review of Drift King, a unique racing game with multiple bonuses, and play for free. Enjoy challenging gameplay and immerse in high-speed visuals.", 2) | Out-Null
